$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be
# auto-converted to numbers by Excel (these are text fields in the sheet).
$ws.Range("D2").Value = '29.376.95'
$ws.Range("E2").Value = '  -0.58%  '
$ws.Range("D3").Value = '1.840.24'
$ws.Range("E3").Value = '  -0.54%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9998'
$ws.Range("E4").Value = '  +0.22%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '239.11'
$ws.Range("E5").Value = '  -0.60%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6265'
$ws.Range("E6").Value = '  -0.66%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.001'
$ws.Range("E7").Value = '  +0.19%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07394'
$ws.Range("E8").Value = '  -1.35%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2888'
$ws.Range("E9").Value = '  -1.01%  '
$ws.Range("E10").Value = '  +0.47%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07722'
$ws.Range("D12").Value = '1.831.51'
$ws.Range("E12").Value = '  -0.95%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.958'
$ws.Range("E13").Value = '  -1.26%  '
$ws.Range("E14").Value = '  -2.25%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.00001037'
$ws.Range("E15").Value = '  -0.91%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '81.52'
$ws.Range("E16").Value = '  -0.96%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.252'
$ws.Range("E17").Value = '  -0.25%  '
$ws.Range("D18").Value = '29.342.70'
$ws.Range("E18").Value = '  -0.64%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '234.16'
$ws.Range("E19").Value = '  +1.74%  '
$ws.Range("E20").Value = '  -1.09%  '
$ws.Range("E21").Value = '  +0.22%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.289'
$ws.Range("E22").Value = '  -3.65%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '157.06'
$ws.Range("E24").Value = '  -1.57%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '8.463'
$ws.Range("E25").Value = '  -0.86%  '
$ws.Range("E26").Value = '  -2.28%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.07271'
$ws.Range("E28").Value = '  +11.19%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.495'
$ws.Range("E29").Value = '  +4.72%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.482'
$ws.Range("E30").Value = '  -0.49%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.025'
$ws.Range("E31").Value = '  -2.00%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.031'
$ws.Range("E32").Value = '  -1.87%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.157'
$ws.Range("E33").Value = '  +0.68%  '
$ws.Range("E34").Value = '  -1.18%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7115'
$ws.Range("E35").Value = '  +1.68%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.580'
$ws.Range("E36").Value = '  +0.12%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01832'
$ws.Range("E37").Value = '  -1.90%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.787'
$ws.Range("E38").Value = '  -1.85%  '
$ws.Range("D39").Value = '1.233.46'
$ws.Range("E39").Value = '  -2.59%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.772'
$ws.Range("E40").Value = '  -1.08%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9506'
$ws.Range("E41").Value = '  +1.12%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.001'
$ws.Range("E42").Value = '  +0.13%  '
$ws.Range("D43").Value = '1.992.32'
$ws.Range("E43").Value = '  -1.59%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '101.10'
$ws.Range("E44").Value = '  -0.24%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '65.17'
$ws.Range("E45").Value = '  -1.83%  '
$ws.Range("E46").Value = '  +0.01%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.693'
$ws.Range("E47").Value = '  -2.92%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '6.945'
$ws.Range("E48").Value = '  -2.40%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.922'
$ws.Range("E49").Value = '  -1.39%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.1129'
$ws.Range("E50").Value = '  -3.43%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.3878'
$ws.Range("E51").Value = '  -2.31%  '
